$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as TEXT (matches original inlineStr cells)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '57.866.04'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '3.136.05'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '529.07'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').Value = '138.83'
$ws.Range('E6').Value = '  -1.39%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.134.76'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('E9').Value = '  +3.23%  '
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('D12').Value = '0.396'
$ws.Range('E12').Value = '  +2.72%  '
$ws.Range('D13').Value = '3.677.13'
$ws.Range('E13').Value = '  +1.02%  '
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = '25.48'
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').Value = '58.003.25'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '3.134.71'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').Value = '6.02'
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('D20').Value = '12.80'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = '7.97'
$ws.Range('E21').Value = '  -1.28%  '
$ws.Range('D22').Value = '354.10'
$ws.Range('E22').Value = '  +5.16%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').Value = '68.83'
$ws.Range('E24').Value = '  +3.36%  '
$ws.Range('E25').Value = '  -0.85%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').Value = '0.995'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').Value = '0.0₃0915'
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('E29').Value = '  +3.88%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').Value = '6.19'
$ws.Range('E31').Value = '  -5.30%  '
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('D33').Value = '21.20'
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('E34').Value = '  -1.03%  '
$ws.Range('E35').Value = '  +7.01%  '
$ws.Range('D36').Value = '158.72'
$ws.Range('E36').Value = '  +1.81%  '
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('D38').Value = '26.62'
$ws.Range('E38').Value = '  -1.45%  '
$ws.Range('E39').Value = '  -1.49%  '
$ws.Range('D40').Value = '0.0671'
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '4.18'
$ws.Range('E41').Value = '  +6.81%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.62'
$ws.Range('E42').Value = '  +5.70%  '
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('D44').Value = '3.171.05'
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('E45').Value = '  +5.13%  '
$ws.Range('D46').Value = '36.55'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('D48').Value = '2.314.79'
$ws.Range('E48').Value = '  +0.55%  '
$ws.Range('D49').Value = '0.970'
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('D50').Value = '20.47'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').Value = '6.04'
$ws.Range('E51').Value = '  +0.43%  '
